# Masculinity_Data_Dictionary.xlsx - "drop NA vals, wrote more explanation markdown"
#
# This script edits the "Cleaning needed" sheet:
#  - replaces the generic "If no answer for all, then drop line?" flag text
#    in column D with 7 distinct, numbered "Create new col to flag no answer
#    for any col" markers - one per logical question-group of rows.
#  - fixes a typo ("secual" -> "sexual") in the harassment question text
#  - adds a new explanatory note about a small-subset-only question
#  - marks several "not interesting / not relevant" rows with strikethrough
#  - tidies up the view (selection / frozen scroll position) and widens
#    column A

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cleaning needed")

# ---------------------------------------------------------------------
# 1) Column D: replace the old blanket "If no answer for all, then drop
#    line?" label with 7 specific group markers.
# ---------------------------------------------------------------------

$group1 = "1- Create new col to flag no answer for any col"
$group2 = "2-Create new col to flag no answer for any col"
$group3 = "3- Create new col to flag no answer for any col"
$group4 = "4- Create new col to flag no answer for any col"
$group5 = "5- Create new col to flag no answer for any col"
$group6 = "6- Create new col to flag no answer for any col"
$group7 = "7- Create new col to flag no answer for any col"

# Group 1: source_ideas_* (rows 6-11)
$ws.Range("D6:D11").Value = $group1

# Group 2: worry_* (rows 24-35)
$ws.Range("D24:D35").Value = $group2

# Group 3: men_earn_more ... men_advantage_other (rows 37-44)
$ws.Range("D37:D44").Value = $group3

# Group 4: hire_women ... men_disadvantage_other (rows 45-49)
$ws.Range("D45:D49").Value = $group4

# ---------------------------------------------------------------------
# 2) Fix the "secual" -> "sexual" typo in the harassment question text,
#    and add the new "small subset" explanatory note.
# ---------------------------------------------------------------------

$ws.Range("E50").Value = "Have you seen or heard of a sexual harassment incident at your work? If so, how did you respond?"
$ws.Range("E57").Value = "Only answered by a small subset of respondents"

# Group 5: harass_confront ... harass_other (rows 50-56) - these rows had
# no value in column D before.
$ws.Range("D50:D56").Value = $group5

# Group 6: int_body_lang ... int_other (rows 69-74) - new D values.
$ws.Range("D69:D74").Value = $group6

# Group 7: bound_wonder ... bound_none (rows 75-78) - new D values.
$ws.Range("D75:D78").Value = $group7

# ---------------------------------------------------------------------
# 3) Strikethrough formatting for the "not interesting / not relevant"
#    rows (device, race2, racethn4, educ3, educ4, age3, kids, orientation,
#    weight) - rows 90-98, columns A:D.
# ---------------------------------------------------------------------

$ws.Range("A90:D98").Font.Strikethrough = $true

# ---------------------------------------------------------------------
# 4) View tidy-up: clear the scrolled "topLeftCell", move the selection
#    to the newly-edited rows, and widen column A so the longer labels
#    are easier to read.
# ---------------------------------------------------------------------

$ws.Range("A1").Select()
$ws.Range("A12:XFD23").Select()
$ws.Columns.Item(1).ColumnWidth = 27.8
